# Add the "Needed to reject it" reason to the two rejected test-case rows
# (TestScenario_4 / Delete Account) in the ReasonToReject column (J), and
# move the active selection from I19 to J19, widening column J to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J9").Value  = "Needed to reject it"
$ws.Range("J14").Value = "Needed to reject it"

# Widen column J (ReasonToReject) so the new text is readable.
$ws.Columns.Item(10).ColumnWidth = 24.65

# Move the selection to J19 (was I19).
$ws.Range("J19").Select() | Out-Null
